$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidential disclaimer text: date changes from 2021-04-08 to 2021-04-09
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update the D and E column values (weights / percent changes) for rows 2-8
$ws.Range("D2").Value = 0.4924531748675571
$ws.Range("E2").Value = 0.005656029606870838

$ws.Range("D3").Value = 0.2505812050810385
$ws.Range("E3").Value = 0.008355321020228734

$ws.Range("D4").Value = 0.09834816426333637
$ws.Range("E4").Value = 0.00387984981226519

$ws.Range("D5").Value = 0.1011284141341417
$ws.Range("E5").Value = 0.005699381761978373

$ws.Range("D6").Value = 0.02962929676090144
$ws.Range("E6").Value = 0.00257171117705246

$ws.Range("D7").Value = 0.02785974489302485
$ws.Range("E7").Value = 0.001234282187765112

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.005947546473449528

# Restore sheet protection (it was protected before this edit)
$ws.Protect()
